# Update the "Metadata" worksheet of the CodeSystem workbook:
#  - Version bumped from 5.0.0 to 6.0.0
#  - Date bumped to the new publish timestamp
#  - Publisher value filled in ("Alvearie Team")
#  - The old duplicate "Contact" / "No display for ContactDetail" row is
#    replaced by a new "Jurisdiction" / "United States of America" row
#  - "Case Sensitive" value filled in ("true")
#  - Everything below that shifts up by one row, shrinking the used range
#    from A1:B22 down to A1:B21

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version
$ws.Range("B3").Value = "6.0.0"

# Date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher gets a value now
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 used to be the second "Contact" row (a duplicate); it becomes the
# new "Jurisdiction" row.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Everything from the old row 12 ("Description") through row 22 ("Count")
# shifts up one row to row 11 through row 21.
$ws.Range("A11").Value = "Description"
$ws.Range("B11").Value = "Codes indicating conversation types for Engagement communications"

$ws.Range("A12").Value = "Purpose"
$ws.Range("B12").Value = ""

$ws.Range("A13").Value = "Copyright"
$ws.Range("B13").Value = ""

$ws.Range("A14").Value = "Case Sensitive"
# Plain "true"/"false" text would be auto-coerced to a boolean by Value's
# smart-input parsing, so force text entry with a leading quote-prefix and
# then restore the original (non-quote-prefixed) cell formatting.
$ws.Range("B14").Value = "'true"
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)

$ws.Range("A15").Value = "Value Set (all codes)"
$ws.Range("B15").Value = ""

$ws.Range("A16").Value = "Hierarchy"
$ws.Range("B16").Value = ""

$ws.Range("A17").Value = "Compositional"
$ws.Range("B17").Value = ""

$ws.Range("A18").Value = "Version Needed?"
$ws.Range("B18").Value = ""

$ws.Range("A19").Value = "Content"
$ws.Range("B19").Value = "complete"

$ws.Range("A20").Value = "Supplements"
$ws.Range("B20").Value = ""

$ws.Range("A21").Value = "Count"
# Likewise, a bare numeric-looking string would become a Number cell, but
# the source data keeps "2" as text, so use the same quote-prefix trick.
$ws.Range("B21").Value = "'2"
$ws.Range("B20").Copy()
$ws.Range("B21").PasteSpecial(-4122)

# The table is now one row shorter: delete the old trailing row 22 so the
# used range becomes A1:B21.
$ws.Rows.Item(22).Delete()
